$d = $word.ActiveDocument

$pairs = @(
    @{ Old = "32×34=1088"; New = "73×47=3431" },
    @{ Old = "70×98=6860"; New = "54×11=594" },
    @{ Old = "24×19=456"; New = "72×67=4824" },
    @{ Old = "14×87=1218"; New = "69×56=3864" },
    @{ Old = "78×54=4212"; New = "40×31=1240" },
    @{ Old = "86×86=7396"; New = "51×88=4488" },
    @{ Old = "11×44=484"; New = "32×41=1312" },
    @{ Old = "70×71=4970"; New = "45×41=1845" },
    @{ Old = "67×79=5293"; New = "24×21=504" },
    @{ Old = "30×69=2070"; New = "49×17=833" },
    @{ Old = "28×89=2492"; New = "97×95=9215" },
    @{ Old = "29×48=1392"; New = "56×39=2184" },
    @{ Old = "99×84=8316"; New = "34×50=1700" },
    @{ Old = "81×71=5751"; New = "83×15=1245" },
    @{ Old = "47×18=846"; New = "47×79=3713" },
    @{ Old = "58×84=4872"; New = "34×29=986" },
    @{ Old = "17×78=1326"; New = "76×78=5928" },
    @{ Old = "21×54=1134"; New = "23×76=1748" },
    @{ Old = "91×74=6734"; New = "30×46=1380" },
    @{ Old = "42×34=1428"; New = "99×92=9108" },
    @{ Old = "66×98=6468"; New = "56×85=4760" },
    @{ Old = "48×64=3072"; New = "74×20=1480" },
    @{ Old = "15×62=930"; New = "80×32=2560" },
    @{ Old = "65×70=4550"; New = "38×74=2812" },
    @{ Old = "64×64=4096"; New = "13×27=351" }
)

foreach ($pair in $pairs) {
    $range = $d.Content
    $range.Find.Execute($pair.Old, $true, $false, $false, $false, $false, $true, 1, $false, $pair.New, 2)
}
